# Update "views/popularity" style counter values (column F) across the
# worksheets of the 北京-漫展信息 workbook. Each entry is keyed by the
# worksheet name (rather than index) so the script is robust to sheet
# ordering.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2  = 306
        3  = 500
        4  = 411
        5  = 8501
        7  = 10593
        22 = 1806
        23 = 69
        24 = 535
        26 = 284
        27 = 59
        30 = 1171
        31 = 22
        34 = 437
        40 = 345
        41 = 92
        42 = 282
        43 = 637
        45 = 94
        46 = 90
    }
    "演出" = @{
        6  = 41
        17 = 381
    }
    "本地生活" = @{
        3 = 2798
    }
    "全部类型" = @{
        2  = 306
        4  = 500
        9  = 411
        10 = 8501
        12 = 10593
        19 = 1806
        20 = 69
        21 = 535
        22 = 284
        23 = 59
        26 = 41
        28 = 1171
        29 = 22
        35 = 437
        41 = 345
        42 = 92
        43 = 282
        46 = 381
        47 = 637
        48 = 94
        49 = 90
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
